$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a blank row at the top; the old header row (r2) shifts down to r3.
$ws.Rows.Item(1).Insert()

# 2. New title row: "THÔNG TIN IMPORT Nguồn Thu" in A1 (bold Times New Roman 11),
#    spanning visually over A1:B1.
$ws.Range("A1").Value = "THÔNG TIN IMPORT Nguồn Thu"
$ws.Range("A1:B1").Font.Name = "Times New Roman"
$ws.Range("A1:B1").Font.Size = 11
$ws.Range("A1:B1").Font.Bold = $true

# 3. Delete the now-empty row 2 so the still-intact old header row (now r3) shifts back
#    up to r2, keeping its row-level height/style attributes and every cell's existing
#    per-cell style intact (no cut/paste needed).
$ws.Rows.Item(2).Delete()

# 4. Column A header text changes from "Tên dự án" to "Tên nguồn thu" (A2 already carries
#    the correct header cell style, since it's the very same cell that held "Tên dự án").
$ws.Range("A2").Value = "Tên nguồn thu"

# 5. Header row fill color: grey tint -> solid orange (FFC000) for the whole header row.
$ws.Range("A2:F2").Interior.Color = 49407   # RGB(255,192,0) = FFC000

# 6. Row height for the header row.
$ws.Rows.Item(2).RowHeight = 31.5

# 7. Column widths (nearest values the host's pixel model can represent).
$ws.Columns.Item(1).ColumnWidth = 33.5
$ws.Columns.Item(2).ColumnWidth = 37.666666666666664
$ws.Columns.Item(3).ColumnWidth = 19.666666666666668
$ws.Columns.Item(4).ColumnWidth = 22
$ws.Columns.Item(5).ColumnWidth = 17.666666666666668
$ws.Columns.Item(6).ColumnWidth = 30.166666666666668

# 8. Selection, matching the saved view in the target file.
[void]$ws.Range("B10").Select()

Write-Output "done"
